# Workbook: evidence.xlsx
# Re-upload with updated evidence values on sheet "A7" and a change of
# which sheet/tab + cell selection is active.

$wb = $excel.ActiveWorkbook

# --- Sheet "A20": update the cell selection (no longer the active tab) ---
$wsA20 = $wb.Worksheets.Item("A20")
$wsA20.Range("C8").Select()

# --- Sheet "A7": fill in new evidence (ibc class on chain / nft id) ---
$wsA7 = $wb.Worksheets.Item("A7")
$wsA7.Range("A2").Value = "ibc/EA5F1D6953610D5BE66EDBA85DF832365D3F775044C24A52C3193755F9C43BF8"
$wsA7.Range("B2").Value = "beauty008"

# Make "A7" the active sheet/tab with its own cell selection.
$wsA7.Activate()
$wsA7.Range("B2").Select()
